$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# text (quote-prefix) so Excel keeps them as strings, matching the
# original inlineStr/text cell type instead of auto-converting to a number.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '67.410.76'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '3.494.58'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").Value = '  +0.01%  '

Set-TextValue $ws.Range("D5") '598.46'
$ws.Range("E5").Value = '  +0.63%  '

Set-TextValue $ws.Range("D6") '179.61'
$ws.Range("E6").Value = '  +3.63%  '

Set-TextValue $ws.Range("D7") '0.609'
$ws.Range("E7").Value = '  +4.12%  '

$ws.Range("D9").Value = '3.496.07'
$ws.Range("E9").Value = '  -0.11%  '

Set-TextValue $ws.Range("D10") '0.138'
$ws.Range("E10").Value = '  +4.40%  '

Set-TextValue $ws.Range("D11") '7.03'
$ws.Range("E11").Value = '  -1.88%  '

Set-TextValue $ws.Range("D12") '0.436'
$ws.Range("E12").Value = '  +1.04%  '

$ws.Range("D13").Value = '4.101.44'
$ws.Range("E13").Value = '  -0.09%  '

Set-TextValue $ws.Range("D14") '32.21'
$ws.Range("E14").Value = '  +9.12%  '

Set-TextValue $ws.Range("D15") '0.135'
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("D16").Value = '67.416.50'
$ws.Range("E16").Value = '  +0.68%  '

Set-TextValue $ws.Range("D17") '0.0000178'
$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("D18").Value = '3.500.07'
$ws.Range("E18").Value = '  +0.12%  '

Set-TextValue $ws.Range("D19") '6.29'
$ws.Range("E19").Value = '  +0.03%  '

Set-TextValue $ws.Range("D20") '14.28'
$ws.Range("E20").Value = '  +0.11%  '

Set-TextValue $ws.Range("D21") '390.61'
$ws.Range("E21").Value = '  -0.68%  '

Set-TextValue $ws.Range("D22") '7.94'
$ws.Range("E22").Value = '  -0.10%  '

Set-TextValue $ws.Range("D23") '73.67'
$ws.Range("E23").Value = '  +0.46%  '

Set-TextValue $ws.Range("D24") '0.541'
$ws.Range("E24").Value = '  +1.11%  '

Set-TextValue $ws.Range("D25") '0.997'
$ws.Range("E25").Value = '  -0.33%  '

Set-TextValue $ws.Range("D26") '5.73'
$ws.Range("E26").Value = '  +0.74%  '

Set-TextValue $ws.Range("D27") '0.0000122'
$ws.Range("E27").Value = '  +0.49%  '

Set-TextValue $ws.Range("D28") '10.35'
$ws.Range("E28").Value = '  +0.85%  '

Set-TextValue $ws.Range("D29") '0.175'
$ws.Range("E29").Value = '  -3.40%  '

$ws.Range("E30").Value = '  +0.54%  '

Set-TextValue $ws.Range("D31") '6.15'
$ws.Range("E31").Value = '  +0.10%  '

Set-TextValue $ws.Range("D32") '1.42'
$ws.Range("E32").Value = '  -0.21%  '

Set-TextValue $ws.Range("D33") '2.07'
$ws.Range("E33").Value = '  +0.76%  '

Set-TextValue $ws.Range("D34") '23.51'
$ws.Range("E34").Value = '  -0.67%  '

Set-TextValue $ws.Range("D35") '7.39'
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  -0.04%  '

Set-TextValue $ws.Range("D37") '1.60'
$ws.Range("E37").Value = '  -0.74%  '

Set-TextValue $ws.Range("D38") '162.70'
$ws.Range("E38").Value = '  -1.13%  '

Set-TextValue $ws.Range("D39") '0.879'
$ws.Range("E39").Value = '  +0.08%  '

Set-TextValue $ws.Range("D40") '2.82'
$ws.Range("E40").Value = '  +10.39%  '

Set-TextValue $ws.Range("D41") '1.88'
$ws.Range("E41").Value = '  -1.30%  '

Set-TextValue $ws.Range("D42") '6.83'
$ws.Range("E42").Value = '  -0.63%  '

Set-TextValue $ws.Range("D43") '4.64'
$ws.Range("E43").Value = '  -0.17%  '

$ws.Range("D44").Value = '2.848.45'
$ws.Range("E44").Value = '  -0.05%  '

Set-TextValue $ws.Range("D45") '26.49'
$ws.Range("E45").Value = '  +1.73%  '

Set-TextValue $ws.Range("D46") '26.69'
$ws.Range("E46").Value = '  -2.81%  '

Set-TextValue $ws.Range("D47") '0.0722'
$ws.Range("E47").Value = '  -2.10%  '

$ws.Range("E48").Value = '  -2.20%  '

Set-TextValue $ws.Range("D49") '0.0300'
$ws.Range("E49").Value = '  -0.55%  '

Set-TextValue $ws.Range("D50") '334.67'
$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("E51").Value = '  -1.70%  '
